$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E15").Value = "DX11Renderer STD::COUT Description"
$ws.Range("C15").Value = "std::cout logging to defined function"

$ws.Range("C15").Select() | Out-Null
